$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 11 values (B, C, D) - A11 stays "10960370"
$ws.Range("B11").Value = "828959809"
$ws.Range("C11").Value = "3016877411"
$ws.Range("D11").Value = "732111198172294"

# Update row 12 values (B, C, D) - A12 stays "10960370"
$ws.Range("B12").Value = "12669894"
$ws.Range("C12").Value = "3016876876"
$ws.Range("D12").Value = "732111198172293"

# Remove row 13 entirely (it was a duplicate of row 11, no longer present)
$ws.Rows.Item(13).Delete()

# Update the selection to match the new active cell used after editing
$ws.Range("C11").Select()
